$wb = $excel.ActiveWorkbook

# Generate Report for Handback:
# The first data row (0c5cca4e... file) on each language sheet gets its
# "Correspond Handoff Datetime" (column E) and "Correspond Handback DateTime"
# (column H) refreshed to new, later timestamps produced by the report run.
# The second data row (d40e0c57... file) keeps its original timestamps.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 13:04:17"
$wsZhCn.Range("H2").Value = "2016-03-24 13:04:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 13:04:21"
$wsDeDe.Range("H2").Value = "2016-03-24 13:04:50"
